# Update gh-pages output (合肥-漫展信息.xlsx) to the state generated at 456a3b4.
#
# Two sheets change:
#   - 展览  (sheet index 1, "Exhibition")  : ticket/price refreshes + one new row (20)
#   - 全部类型 (sheet index 4, "All types") : same refreshes + the same new row
#     inserted before the existing last row (which shifts from 21 -> 22)

$wb = $excel.ActiveWorkbook

function Update-CommonCounts($ws) {
    # "想去人数" (interest count) / price refreshes shared by both sheets
    $ws.Range("F2").Value  = 8820
    $ws.Range("F3").Value  = 8181
    $ws.Range("F9").Value  = 147
    $ws.Range("F10").Value = 200
    $ws.Range("F13").Value = 200
    $ws.Range("F14").Value = 4497
    $ws.Range("G14").Value = 29.9
    $ws.Range("F16").Value = 80
    $ws.Range("F19").Value = 151
    $ws.Range("F20").Value = 128
}

function Set-TextValue($range, [string]$text) {
    # Force a literal text value (avoids Excel auto-converting things that
    # look like dates, e.g. "2024-07-27", into date serials).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet "展览" (index 1): refresh counts, then append the new row 21.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
Update-CommonCounts $ws1

$ws1.Range("A20").Copy($ws1.Range("A21"))
$ws1.Range("A21").Value = 20
Set-TextValue $ws1.Range("B21") "2024-07-27"
$ws1.Range("C21").Value = "安徽·MAX特摄only展"
$ws1.Range("D21").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws1.Range("E21").Value = "2024.07.27 09:30-07.27 18:00"
$ws1.Range("F21").Value = 1
$ws1.Range("G21").Value = 50
$ws1.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=83684"
$ws1.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202404/jv1CiqqW1712029200830.jpeg"

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4): refresh counts, then insert the new row at
# 21 (pushing the existing last row down to 22) and fill both rows.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Update-CommonCounts $ws4

# Push the current row 21 (合肥·首届包河留声机音乐节...) down to row 22,
# carrying its values + formatting intact, then renumber its index.
$ws4.Range("A21:I21").Copy($ws4.Range("A22"))
$ws4.Range("A22").Value = 21

# Overwrite row 21 in place with the new record (keeps the existing A21
# "index column" style, which already matches the other rows).
$ws4.Range("A21").Value = 20
Set-TextValue $ws4.Range("B21") "2024-07-27"
$ws4.Range("C21").Value = "安徽·MAX特摄only展"
$ws4.Range("D21").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws4.Range("E21").Value = "2024.07.27 09:30-07.27 18:00"
$ws4.Range("F21").Value = 1
$ws4.Range("G21").Value = 50
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=83684"
$ws4.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202404/jv1CiqqW1712029200830.jpeg"
